# Update "想去人数" (number of people interested) figures for several
# events, as reflected in the regenerated gh-pages data output.
#
# Sheet "展览" (Exhibitions): rows 3-8, column F
# Sheet "全部类型" (All types): rows 3-6 and 8-9, column F
# (sheet "演出" and "本地生活" are untouched)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 192
$wsExpo.Range("F4").Value = 2224
$wsExpo.Range("F5").Value = 1708
$wsExpo.Range("F6").Value = 321
$wsExpo.Range("F7").Value = 95
$wsExpo.Range("F8").Value = 766

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 192
$wsAll.Range("F4").Value = 2224
$wsAll.Range("F5").Value = 1708
$wsAll.Range("F6").Value = 321
$wsAll.Range("F8").Value = 95
$wsAll.Range("F9").Value = 766
